# Apply corrections to the District (column G) names, using the official
# name "Belagavi (Belgaum)" for every district row that is not "Raichur".
# Also clear the stray empty inline-string cell at F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44
for ($r = 4; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Text -ne "Raichur") {
        $cell.Value = "Belagavi (Belgaum)"
    }
}

# Remove the empty stray cell F3 entirely.
$ws.Cells.Item(3, 6).ClearContents()
